$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3625
$ws.Range("I62").Value = 3625
$ws.Range("K62").Value = 3625
$ws.Range("M62").Value = -3001
$ws.Range("H65").Value = 3625
$ws.Range("I65").Value = 3625
$ws.Range("K65").Value = 18125
$ws.Range("M65").Value = -15005
$ws.Range("H80").Value = 1151.2858
$ws.Range("I80").Value = 917
$ws.Range("J80").Value = 1737
$ws.Range("K80").Value = 2751
$ws.Range("L80").Value = 5211
$ws.Range("M80").Value = -1753
$ws.Range("N80").Value = -7207
$ws.Range("H83").Value = 1151.2858
$ws.Range("I83").Value = 917
$ws.Range("J83").Value = 1737
$ws.Range("K83").Value = 8253
$ws.Range("L83").Value = 15633
$ws.Range("M83").Value = -3261
$ws.Range("N83").Value = -25617
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H103").Value = 1798.5714
$ws.Range("I103").Value = 1838
$ws.Range("K103").Value = 5514
$ws.Range("M103").Value = -4928
$ws.Range("H112").Value = 2599.9443
$ws.Range("J112").Value = 3065.3845
$ws.Range("L112").Value = 9196.1535
$ws.Range("N112").Value = -11412.1535
$ws.Range("H135").Value = 660.3333
$ws.Range("I135").Value = 605.375
$ws.Range("K135").Value = 5448.375
$ws.Range("M135").Value = -2913.375
$ws.Range("H137").Value = 2695
$ws.Range("I137").Value = 2367.25
$ws.Range("J137").Value = 4006
$ws.Range("K137").Value = 7101.75
$ws.Range("L137").Value = 12018
$ws.Range("M137").Value = -4551.75
$ws.Range("N137").Value = -17118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 2313.5334
$ws.Range("I32").Value = 2134.6897
$ws.Range("K32").Value = 2134.6897
$ws.Range("M32").Value = -1847.6897
$ws.Range("H38").Value = 14999.5
$ws.Range("I38").Value = 14999.5
$ws.Range("K38").Value = 14999.5
$ws.Range("M38").Value = -14532.5
$ws.Range("H45").Value = 2064.1667
$ws.Range("J45").Value = 2266
$ws.Range("L45").Value = 2266
$ws.Range("N45").Value = -3020
$ws.Range("H74").Value = 1060.7646
$ws.Range("I74").Value = 1058.3125
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 1058.3125
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -184.3125
$ws.Range("N74").Value = -2848
$ws.Range("H77").Value = 1060.7646
$ws.Range("I77").Value = 1058.3125
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 5291.5625
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -923.5625
$ws.Range("N77").Value = -14236
$ws.Range("H92").Value = 34997.5
$ws.Range("J92").Value = 34997.5
$ws.Range("L92").Value = 34997.5
$ws.Range("N92").Value = -39989.5
$ws.Range("H96").Value = 12000
$ws.Range("J96").Value = 12000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -17492
$ws.Range("H132").Value = 3398
$ws.Range("I132").Value = 3398
$ws.Range("K132").Value = 10194
$ws.Range("M132").Value = -7664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2259.1428
$ws.Range("I94").Value = 2481.8
$ws.Range("K94").Value = 2481.8
$ws.Range("M94").Value = -2030.8
$ws.Range("H99").Value = 4569.5
$ws.Range("I99").Value = 5177.857
$ws.Range("J99").Value = 3150
$ws.Range("K99").Value = 5177.857
$ws.Range("L99").Value = 3150
$ws.Range("M99").Value = -3679.857
$ws.Range("N99").Value = -6146
$ws.Range("H134").Value = 9053.357
$ws.Range("I134").Value = 9442.076999999999
$ws.Range("K134").Value = 28326.231
$ws.Range("M134").Value = -25791.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 6113.5713
$ws.Range("J41").Value = 21250
$ws.Range("L41").Value = 21250
$ws.Range("N41").Value = -22106
$ws.Range("H99").Value = 2122.111
$ws.Range("I99").Value = 2975.2
$ws.Range("J99").Value = 1055.75
$ws.Range("K99").Value = 2975.2
$ws.Range("L99").Value = 1055.75
$ws.Range("M99").Value = -1477.2
$ws.Range("N99").Value = -4051.75
$ws.Range("H126").Value = 2122.111
$ws.Range("I126").Value = 2975.2
$ws.Range("J126").Value = 1055.75
$ws.Range("K126").Value = 8925.599999999999
$ws.Range("L126").Value = 3167.25
$ws.Range("M126").Value = -6455.599999999999
$ws.Range("N126").Value = -8107.25
$ws.Range("H134").Value = 1999.5
$ws.Range("I134").Value = 1999.5
$ws.Range("K134").Value = 5998.5
$ws.Range("M134").Value = -3463.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 12999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 12999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 116991
$ws.Range("N132").Value = -122051
$ws.Range("M132").ClearContents()
$ws.Range("H134").Value = 457.4
$ws.Range("I134").Value = 457.4
$ws.Range("K134").Value = 1372.2
$ws.Range("M134").Value = 3697.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 15250
$ws.Range("J43").Value = 15250
$ws.Range("L43").Value = 15250
$ws.Range("N43").Value = -15552
$ws.Range("H80").Value = 1654.3334
$ws.Range("I80").Value = 1125.25
$ws.Range("K80").Value = 1125.25
$ws.Range("M80").Value = -127.25
$ws.Range("H83").Value = 1654.3334
$ws.Range("I83").Value = 1125.25
$ws.Range("K83").Value = 5626.25
$ws.Range("M83").Value = -634.25
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1099.1538
$ws.Range("I22").Value = 935.5454999999999
$ws.Range("J22").Value = 1999
$ws.Range("K22").Value = 935.5454999999999
$ws.Range("L22").Value = 1999
$ws.Range("M22").Value = -640.5454999999999
$ws.Range("N22").Value = -2589
$ws.Range("H27").Value = 1099.1538
$ws.Range("I27").Value = 935.5454999999999
$ws.Range("J27").Value = 1999
$ws.Range("K27").Value = 935.5454999999999
$ws.Range("L27").Value = 1999
$ws.Range("M27").Value = -828.5454999999999
$ws.Range("N27").Value = -2213
$ws.Range("H46").Value = 2792.2424
$ws.Range("J46").Value = 3302.647
$ws.Range("L46").Value = 3302.647
$ws.Range("N46").Value = -3678.647
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H55").Value = 173.82608
$ws.Range("I55").Value = 139.35715
$ws.Range("J55").Value = 227.44444
$ws.Range("K55").Value = 139.35715
$ws.Range("L55").Value = 227.44444
$ws.Range("M55").Value = 33.64285000000001
$ws.Range("N55").Value = -573.44444
$ws.Range("H68").Value = 52747.25
$ws.Range("I68").Value = 3991
$ws.Range("J68").Value = 68999.336
$ws.Range("K68").Value = 3991
$ws.Range("L68").Value = 68999.336
$ws.Range("M68").Value = -3242
$ws.Range("N68").Value = -70497.336
$ws.Range("H71").Value = 52747.25
$ws.Range("I71").Value = 3991
$ws.Range("J71").Value = 68999.336
$ws.Range("K71").Value = 19955
$ws.Range("L71").Value = 344996.68
$ws.Range("M71").Value = -16211
$ws.Range("N71").Value = -352484.68
$ws.Range("H82").Value = 1655.2858
$ws.Range("I82").Value = 1357.4
$ws.Range("J82").Value = 2400
$ws.Range("K82").Value = 1357.4
$ws.Range("L82").Value = 2400
$ws.Range("M82").Value = -996.4000000000001
$ws.Range("N82").Value = -3122
$ws.Range("H85").Value = 1655.2858
$ws.Range("I85").Value = 1357.4
$ws.Range("J85").Value = 2400
$ws.Range("K85").Value = 1357.4
$ws.Range("L85").Value = 2400
$ws.Range("M85").Value = -109.4000000000001
$ws.Range("N85").Value = -4896
$ws.Range("H132").Value = 3056.4443
$ws.Range("I132").Value = 2462.7693
$ws.Range("K132").Value = 7388.3079
$ws.Range("M132").Value = -4858.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2100
$ws.Range("I14").Value = 2100
$ws.Range("K14").Value = 2100
$ws.Range("M14").Value = -1932
$ws.Range("H81").Value = 2828.077
$ws.Range("I81").Value = 2897.9167
$ws.Range("J81").Value = 1990
$ws.Range("K81").Value = 5795.8334
$ws.Range("L81").Value = 3980
$ws.Range("M81").Value = -4734.8334
$ws.Range("N81").Value = -6102
$ws.Range("H84").Value = 2828.077
$ws.Range("I84").Value = 2897.9167
$ws.Range("J84").Value = 1990
$ws.Range("K84").Value = 28979.167
$ws.Range("L84").Value = 19900
$ws.Range("M84").Value = -23675.167
$ws.Range("N84").Value = -30508
$ws.Range("H126").Value = 1553.125
$ws.Range("I126").Value = 1385.2
$ws.Range("K126").Value = 4155.6
$ws.Range("M126").Value = -1685.6
$ws.Range("H132").Value = 1377.25
$ws.Range("I132").Value = 1301.3334
$ws.Range("J132").Value = 1605
$ws.Range("K132").Value = 3904.0002
$ws.Range("L132").Value = 4815
$ws.Range("M132").Value = -1374.0002
$ws.Range("N132").Value = -9875
